$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 7.5047619047619
$ws.Range("K2").Value = 10.5335
$ws.Range("AA2").Value = 4.1206818
$ws.Range("BQ2").Value = 157.6
$ws.Range("CC2").Value = 92
$ws.Range("E3").Value = 7.0105263157895
$ws.Range("K3").Value = 7.426
$ws.Range("AA3").Value = 3.08032525
$ws.Range("AP3").Value = 6
$ws.Range("AR3").Value = 23
$ws.Range("AS3").Value = 7
$ws.Range("BA3").Value = 29
$ws.Range("BB3").Value = 37.662337662338
$ws.Range("BC3").Value = 22
$ws.Range("BD3").Value = 37.931034482759
$ws.Range("BG3").Value = 103
$ws.Range("BQ3").Value = 133.2
$ws.Range("CU3").Value = 48
$ws.Range("DA3").Value = 4
$ws.Range("DB3").Value = 66.666666666667
$ws.Range("AA4").Value = 1.09891517
$ws.Range("AQ4").Value = 6
$ws.Range("AR4").Value = 25
$ws.Range("AA7").Value = 1.73291931
$ws.Range("CC7").Value = 35
$ws.Range("K8").Value = 0.4445
$ws.Range("AA8").Value = 0.2298759
$ws.Range("AB8").Value = 216
$ws.Range("AL8").Value = 16.666666666667
$ws.Range("AR8").Value = 15
$ws.Range("AY8").Value = 4
$ws.Range("AZ8").Value = 40
$ws.Range("BA8").Value = 21
$ws.Range("BB8").Value = 52.5
$ws.Range("BC8").Value = 21
$ws.Range("BD8").Value = 53.846153846154
$ws.Range("BG8").Value = 35
$ws.Range("CB8").Value = 5
$ws.Range("CC8").Value = 5
$ws.Range("CS8").Value = 10
$ws.Range("CU8").Value = 19
$ws.Range("CZ8").Value = 6
$ws.Range("AA9").Value = 0.73377377
$ws.Range("AB9").Value = 831
$ws.Range("AA10").Value = 1.34873883
$ws.Range("AB10").Value = 437
$ws.Range("AE10").Value = 310
$ws.Range("AF10").Value = 90.116279069767
$ws.Range("AG10").Value = 344
$ws.Range("AI10").Value = 193
$ws.Range("AO10").Value = 5
$ws.Range("DG10").Value = 217
$ws.Range("E12").Value = 7.205
$ws.Range("AA12").Value = 1.79039029
$ws.Range("AB12").Value = 1688
$ws.Range("AH12").Value = 590
$ws.Range("AI12").Value = 749
$ws.Range("AJ12").Value = 331
$ws.Range("AO12").Value = 16
$ws.Range("AR12").Value = 97
$ws.Range("AV12").Value = 3
$ws.Range("BB12").Value = 54.676258992806
$ws.Range("BD12").Value = 53.508771929825
$ws.Range("BQ12").Value = 144.1
$ws.Range("CB12").Value = 10
$ws.Range("CU12").Value = 63
$ws.Range("DA12").Value = 12
$ws.Range("DB12").Value = 63.157894736842
$ws.Range("DF12").Value = 613
$ws.Range("DG12").Value = 836
$ws.Range("AA15").Value = 0.27437439
$ws.Range("CC15").Value = 32
$ws.Range("AA16").Value = 0.68156774
$ws.Range("CC16").Value = 33
$ws.Range("E18").Value = 6.5625
$ws.Range("AA18").Value = 0.1697539
$ws.Range("BQ18").Value = 52.5
$ws.Range("E19").Value = 7.0176470588235
$ws.Range("AA19").Value = 2.92055973
$ws.Range("AF19").Value = 89.848308051342
$ws.Range("AG19").Value = 857
$ws.Range("AL19").Value = 40.540540540541
$ws.Range("AT19").Value = 36
$ws.Range("BG19").Value = 183
$ws.Range("BQ19").Value = 119.3
$ws.Range("BT19").Value = 87
$ws.Range("CC19").Value = 45
$ws.Range("CD19").Value = 23
$ws.Range("CZ19").Value = 37
$ws.Range("DG19").Value = 493
$ws.Range("E21").Value = 7.3
$ws.Range("AB21").Value = 805
$ws.Range("BQ21").Value = 131.4
$ws.Range("CC21").Value = 66
$ws.Range("CO21").Value = 34
$ws.Range("DL21").Value = 4.2575
